$d = $word.ActiveDocument

# 1. Fix the double space before "adauge" in the third paragraph
# ("... va putea sa  adauge ..." -> "... va putea sa adauge ...")
$d.Content.Find.Execute("sa  adauge", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "sa adauge", 2) | Out-Null

# 2. Append the new paragraphs describing the project implementation.
$end = $d.Content
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "Pentru a realiza functionalitatile descrise mai sus voi crea un proiect nou in InteliiJ IDEA folosind spring framework:")
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "voi realiza conexiunea la baza de date cu ajutorul layer-ului model, in baza de date voi stoca date despre utilizator(date logare, date privind metoda de plata)")
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "voi crea in clasa Controller metode care vor returna un string(nume.html), in care voi realiza felul cum va fi afisat in browser datele. ")
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "Pentru ca utilizatorul sa poata viziona un film vom folosi metoda GET, deoarece utilizatorul trebuie sa primeasca niste data.")
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "Pentru crearea unul cont nou sau logare voi folosi metoda POST unde in body voi trimite catre server datele introduse de utilizator.")
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "Pentru cautarea unui anumit film(dupa denumire) voi folosi la fel metoda GET, unde ca parametru ii voi da numele filmului.")
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "Pentru adaugarea de comentariu voi folosi metoda POST pentru a transmite datele catre server fiind ulterior salvate intr-o baza de date.")
$end.Collapse(0) | Out-Null

$end.InsertAfter("`r" + "Pentru adaugarea metodei de plata la fel voi folosi metoda POST.")
$end.Collapse(0) | Out-Null

# Trailing blank paragraph (InsertBefore on a range collapsed to the very end
# produces a bare "<w:p/>" instead of a paragraph holding an empty run).
$tail = $d.Content
$tail.Collapse(0) | Out-Null
$tail.InsertBefore("`r")

# 3. Turn the two "voi ..." paragraphs into a bulleted list (List Paragraph style).
$listStart = $d.Paragraphs(5).Range.Start
$listEnd = $d.Paragraphs(6).Range.End
$listRange = $d.Range($listStart, $listEnd)
$listRange.Style = "List Paragraph"
$listRange.ListFormat.ApplyBulletDefault()

# 4. Tidy up the generated "List Paragraph" style so it matches Word's usual
#    "Define New Bullet" output (left-indent + no extra space between items).
$listStyle = $d.Styles.Item("List Paragraph")
$listStyle.Priority = 34
$listStyle.ParagraphFormat.LeftIndent = 36
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true
